$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells we touch keep their original text
# representation (no auto-conversion to numbers, so formatting such as
# trailing zeros and thousands separators written as literal dots is
# preserved exactly like the source inline strings).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.852.09"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.120.22"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.09"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.30"
$ws.Range("E6").Value = "  +1.94%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.115.16"
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("E10").Value = "  +11.29%  "

$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("E13").Value = "  +3.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.03"
$ws.Range("E14").Value = "  +4.03%  "

$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.638.95"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.769.44"
$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.15"
$ws.Range("E18").Value = "  -1.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.119.72"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.45"
$ws.Range("E20").Value = "  +2.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.30"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("E24").Value = "  -3.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.06"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.95"
$ws.Range("E27").Value = "  +8.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.04"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("E33").Value = "  -2.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0878"
$ws.Range("E34").Value = "  +9.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("E35").Value = "  +5.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +0.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +13.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.08"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.99"
$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.39"
$ws.Range("E40").Value = "  +6.16%  "

$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0372"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.888.77"
$ws.Range("E43").Value = "  -2.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.72"
$ws.Range("E47").Value = "  +2.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.56"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.69"
$ws.Range("E51").Value = "  -0.77%  "
